$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: "last updated" banner ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 03:22"

# --- Rows 3-221: header + full country table ---
$data = New-Object 'object[,]' 219,8
$data[0,0] = "País"
$data[0,1] = "Casos totales"
$data[0,2] = "Nuevos casos"
$data[0,3] = "Casos activos"
$data[0,4] = "Recuperados"
$data[0,5] = "Casos críticos"
$data[0,6] = "Muertes hoy"
$data[0,7] = "Muertes"
$data[1,0] = "Estados Unidos"
$data[1,1] = 8288278
$data[1,2] = 71687
$data[1,3] = 5395401
$data[1,4] = 2669233
$data[1,5] = 0
$data[1,6] = 928
$data[1,7] = 223644
$data[2,0] = "India"
$data[2,1] = 7430635
$data[2,2] = 65126
$data[2,3] = 6521634
$data[2,4] = 795969
$data[2,5] = 0
$data[2,6] = 886
$data[2,7] = 113032
$data[3,0] = "Brasil"
$data[3,1] = 5201570
$data[3,2] = 30574
$data[3,3] = 4619560
$data[3,4] = 428781
$data[3,5] = 0
$data[3,6] = 716
$data[3,7] = 153229
$data[4,0] = "Rusia"
$data[4,1] = 1369313
$data[4,2] = 15150
$data[4,3] = 1056582
$data[4,4] = 289008
$data[4,5] = 0
$data[4,6] = 232
$data[4,7] = 23723
$data[5,0] = "España"
$data[5,1] = 982723
$data[5,2] = 12169
$data[5,3] = 0
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 222
$data[5,7] = 33775
$data[6,0] = "Argentina"
$data[6,1] = 965609
$data[6,2] = 16546
$data[6,3] = 778501
$data[6,4] = 161385
$data[6,5] = 0
$data[6,6] = 381
$data[6,7] = 25723
$data[7,0] = "Colombia"
$data[7,1] = 945354
$data[7,2] = 8372
$data[7,3] = 837001
$data[7,4] = 79737
$data[7,5] = 0
$data[7,6] = 159
$data[7,7] = 28616
$data[8,0] = "Peru"
$data[8,1] = 862417
$data[8,2] = 2677
$data[8,3] = 769077
$data[8,4] = 59692
$data[8,5] = 0
$data[8,6] = 71
$data[8,7] = 33648
$data[9,0] = "Mexico"
$data[9,1] = 834910
$data[9,2] = 5514
$data[9,3] = 608188
$data[9,4] = 141437
$data[9,5] = 0
$data[9,6] = 387
$data[9,7] = 85285
$data[10,0] = "Francia"
$data[10,1] = 834770
$data[10,2] = 25086
$data[10,3] = 104696
$data[10,4] = 696771
$data[10,5] = 0
$data[10,6] = 178
$data[10,7] = 33303
$data[11,0] = "Sudafrica"
$data[11,1] = 700203
$data[11,2] = 2019
$data[11,3] = 629260
$data[11,4] = 52573
$data[11,5] = 0
$data[11,6] = 61
$data[11,7] = 18370
$data[12,0] = "Reino Unido"
$data[12,1] = 689257
$data[12,2] = 15650
$data[12,3] = 0
$data[12,4] = 0
$data[12,5] = 0
$data[12,6] = 136
$data[12,7] = 43429
$data[13,0] = "Iran"
$data[13,1] = 522387
$data[13,2] = 4552
$data[13,3] = 420910
$data[13,4] = 71607
$data[13,5] = 0
$data[13,6] = 265
$data[13,7] = 29870
$data[14,0] = "Chile"
$data[14,1] = 488190
$data[14,2] = 1694
$data[14,3] = 461097
$data[14,4] = 13564
$data[14,5] = 0
$data[14,6] = 95
$data[14,7] = 13529
$data[15,0] = "Irak"
$data[15,1] = 420303
$data[15,2] = 3501
$data[15,3] = 353962
$data[15,4] = 56199
$data[15,5] = 0
$data[15,6] = 56
$data[15,7] = 10142
$data[16,0] = "Italia"
$data[16,1] = 391611
$data[16,2] = 10010
$data[16,3] = 247872
$data[16,4] = 107312
$data[16,5] = 0
$data[16,6] = 55
$data[16,7] = 36427
$data[17,0] = "Banglades"
$data[17,1] = 386086
$data[17,2] = 1527
$data[17,3] = 300738
$data[17,4] = 79725
$data[17,5] = 0
$data[17,6] = 15
$data[17,7] = 5623
$data[18,0] = "Alemania"
$data[18,1] = 356780
$data[18,2] = 7964
$data[18,3] = 287600
$data[18,4] = 59344
$data[18,5] = 0
$data[18,6] = 26
$data[18,7] = 9836
$data[19,0] = "Indonesia"
$data[19,1] = 353461
$data[19,2] = 4301
$data[19,3] = 277544
$data[19,4] = 63570
$data[19,5] = 0
$data[19,6] = 79
$data[19,7] = 12347
$data[20,0] = "Filipinas"
$data[20,1] = 351750
$data[20,2] = 3139
$data[20,3] = 294865
$data[20,4] = 50354
$data[20,5] = 0
$data[20,6] = 34
$data[20,7] = 6531
$data[21,0] = "Turquia"
$data[21,1] = 343955
$data[21,2] = 1812
$data[21,3] = 301098
$data[21,4] = 33704
$data[21,5] = 0
$data[21,6] = 73
$data[21,7] = 9153
$data[22,0] = "Arabia Saudita"
$data[22,1] = 341495
$data[22,2] = 433
$data[22,3] = 327795
$data[22,4] = 8556
$data[22,5] = 0
$data[22,6] = 17
$data[22,7] = 5144
$data[23,0] = "Pakistan"
$data[23,1] = 321877
$data[23,2] = 659
$data[23,3] = 305835
$data[23,4] = 9421
$data[23,5] = 0
$data[23,6] = 7
$data[23,7] = 6621
$data[24,0] = "Israel"
$data[24,1] = 301896
$data[24,2] = 1695
$data[24,3] = 262503
$data[24,4] = 37252
$data[24,5] = 0
$data[24,6] = 14
$data[24,7] = 2141
$data[25,0] = "Ucrania"
$data[25,1] = 287231
$data[25,2] = 5992
$data[25,3] = 121919
$data[25,4] = 159904
$data[25,5] = 0
$data[25,6] = 106
$data[25,7] = 5408
$data[26,0] = "Paises Bajos"
$data[26,1] = 211938
$data[26,2] = 7984
$data[26,3] = 0
$data[26,4] = 0
$data[26,5] = 0
$data[26,6] = 16
$data[26,7] = 6708
$data[27,0] = "Canada"
$data[27,1] = 194106
$data[27,2] = 2374
$data[27,3] = 163644
$data[27,4] = 20740
$data[27,5] = 0
$data[27,6] = 23
$data[27,7] = 9722
$data[28,0] = "Belgica"
$data[28,1] = 191959
$data[28,2] = 10448
$data[28,3] = 20720
$data[28,4] = 160912
$data[28,5] = 0
$data[28,6] = 49
$data[28,7] = 10327
$data[29,0] = "Rumania"
$data[29,1] = 172516
$data[29,2] = 4026
$data[29,3] = 127076
$data[29,4] = 39691
$data[29,5] = 0
$data[29,6] = 75
$data[29,7] = 5749
$data[30,0] = "Marruecos"
$data[30,1] = 167148
$data[30,2] = 3498
$data[30,3] = 138989
$data[30,4] = 25341
$data[30,5] = 0
$data[30,6] = 46
$data[30,7] = 2818
$data[31,0] = "Chequia"
$data[31,1] = 160112
$data[31,2] = 11102
$data[31,3] = 66093
$data[31,4] = 92736
$data[31,5] = 0
$data[31,6] = 53
$data[31,7] = 1283
$data[32,0] = "Polonia"
$data[32,1] = 157608
$data[32,2] = 7705
$data[32,3] = 87773
$data[32,4] = 66395
$data[32,5] = 0
$data[32,6] = 132
$data[32,7] = 3440
$data[33,0] = "Ecuador"
$data[33,1] = 151659
$data[33,2] = 1299
$data[33,3] = 128134
$data[33,4] = 11168
$data[33,5] = 0
$data[33,6] = 51
$data[33,7] = 12357
$data[34,0] = "Bolivia"
$data[34,1] = 139319
$data[34,2] = 178
$data[34,3] = 103499
$data[34,4] = 27413
$data[34,5] = 0
$data[34,6] = 30
$data[34,7] = 8407
$data[35,0] = "Catar"
$data[35,1] = 128992
$data[35,2] = 189
$data[35,3] = 126006
$data[35,4] = 2764
$data[35,5] = 0
$data[35,6] = 0
$data[35,7] = 222
$data[36,0] = "Nepal"
$data[36,1] = 126137
$data[36,2] = 4392
$data[36,3] = 88040
$data[36,4] = 37382
$data[36,5] = 0
$data[36,6] = 21
$data[36,7] = 715
$data[37,0] = "Panama"
$data[37,1] = 123498
$data[37,2] = 615
$data[37,3] = 99286
$data[37,4] = 21666
$data[37,5] = 0
$data[37,6] = 17
$data[37,7] = 2546
$data[38,0] = "Republica Dominicana"
$data[38,1] = 120450
$data[38,2] = 384
$data[38,3] = 96883
$data[38,4] = 21375
$data[38,5] = 0
$data[38,6] = 3
$data[38,7] = 2192
$data[39,0] = "Kuwait"
$data[39,1] = 114744
$data[39,2] = 729
$data[39,3] = 106495
$data[39,4] = 7559
$data[39,5] = 0
$data[39,6] = 6
$data[39,7] = 690
$data[40,0] = "Emiratos Arabes Unidos"
$data[40,1] = 112849
$data[40,2] = 1412
$data[40,3] = 104943
$data[40,4] = 7451
$data[40,5] = 0
$data[40,6] = 3
$data[40,7] = 455
$data[41,0] = "Kazajistan"
$data[41,1] = 109202
$data[41,2] = 108
$data[41,3] = 104801
$data[41,4] = 2633
$data[41,5] = 0
$data[41,6] = 0
$data[41,7] = 1768
$data[42,0] = "Oman"
$data[42,1] = 108296
$data[42,2] = 0
$data[42,3] = 94229
$data[42,4] = 12996
$data[42,5] = 0
$data[42,6] = 0
$data[42,7] = 1071
$data[43,0] = "Egipto"
$data[43,1] = 105159
$data[43,2] = 126
$data[43,3] = 98089
$data[43,4] = 971
$data[43,5] = 0
$data[43,6] = 11
$data[43,7] = 6099
$data[44,0] = "Suecia"
$data[44,1] = 103200
$data[44,2] = 0
$data[44,3] = 0
$data[44,4] = 0
$data[44,5] = 0
$data[44,6] = 6
$data[44,7] = 5918
$data[45,0] = "Guatemala"
$data[45,1] = 100431
$data[45,2] = 666
$data[45,3] = 89494
$data[45,4] = 7459
$data[45,5] = 0
$data[45,6] = 25
$data[45,7] = 3478
$data[46,0] = "Portugal"
$data[46,1] = 95902
$data[46,2] = 2608
$data[46,3] = 56066
$data[46,4] = 37687
$data[46,5] = 0
$data[46,6] = 21
$data[46,7] = 2149
$data[47,0] = "Costa Rica"
$data[47,1] = 94348
$data[47,2] = 1196
$data[47,3] = 58269
$data[47,4] = 34911
$data[47,5] = 0
$data[47,6] = 9
$data[47,7] = 1168
$data[48,0] = "Japon"
$data[48,1] = 91431
$data[48,2] = 721
$data[48,3] = 84451
$data[48,4] = 5330
$data[48,5] = 0
$data[48,6] = 4
$data[48,7] = 1650
$data[49,0] = "Etiopia"
$data[49,1] = 87834
$data[49,2] = 665
$data[49,3] = 41628
$data[49,4] = 44869
$data[49,5] = 0
$data[49,6] = 12
$data[49,7] = 1337
$data[50,0] = "Bielorrusia"
$data[50,1] = 86392
$data[50,2] = 658
$data[50,3] = 78990
$data[50,4] = 6481
$data[50,5] = 0
$data[50,6] = 5
$data[50,7] = 921
$data[51,0] = "Honduras"
$data[51,1] = 86089
$data[51,2] = 631
$data[51,3] = 34099
$data[51,4] = 49438
$data[51,5] = 0
$data[51,6] = 19
$data[51,7] = 2552
$data[52,0] = "China"
$data[52,1] = 85646
$data[52,2] = 24
$data[52,3] = 80759
$data[52,4] = 253
$data[52,5] = 0
$data[52,6] = 0
$data[52,7] = 4634
$data[53,0] = "Venezuela"
$data[53,1] = 85469
$data[53,2] = 0
$data[53,3] = 77689
$data[53,4] = 7060
$data[53,5] = 0
$data[53,6] = 0
$data[53,7] = 720
$data[54,0] = "Barein"
$data[54,1] = 77325
$data[54,2] = 371
$data[54,3] = 73421
$data[54,4] = 3612
$data[54,5] = 0
$data[54,6] = 3
$data[54,7] = 292
$data[55,0] = "Suiza"
$data[55,1] = 74422
$data[55,2] = 3105
$data[55,3] = 50500
$data[55,4] = 21800
$data[55,5] = 0
$data[55,6] = 7
$data[55,7] = 2122
$data[56,0] = "Moldavia"
$data[56,1] = 65860
$data[56,2] = 784
$data[56,3] = 46543
$data[56,4] = 17768
$data[56,5] = 0
$data[56,6] = 19
$data[56,7] = 1549
$data[57,0] = "Uzbekistan"
$data[57,1] = 62588
$data[57,2] = 310
$data[57,3] = 59624
$data[57,4] = 2444
$data[57,5] = 0
$data[57,6] = 2
$data[57,7] = 520
$data[58,0] = "Armenia"
$data[58,1] = 61460
$data[58,2] = 1465
$data[58,3] = 47541
$data[58,4] = 12863
$data[58,5] = 0
$data[58,6] = 10
$data[58,7] = 1056
$data[59,0] = "Austria"
$data[59,1] = 61387
$data[59,2] = 1163
$data[59,3] = 47618
$data[59,4] = 12887
$data[59,5] = 0
$data[59,6] = 5
$data[59,7] = 882
$data[60,0] = "Nigeria"
$data[60,1] = 61194
$data[60,2] = 212
$data[60,3] = 52304
$data[60,4] = 7771
$data[60,5] = 0
$data[60,6] = 3
$data[60,7] = 1119
$data[61,0] = "Libano"
$data[61,1] = 60113
$data[61,2] = 1368
$data[61,3] = 26468
$data[61,4] = 33136
$data[61,5] = 0
$data[61,6] = 8
$data[61,7] = 509
$data[62,0] = "Singapur"
$data[62,1] = 57901
$data[62,2] = 9
$data[62,3] = 57784
$data[62,4] = 89
$data[62,5] = 0
$data[62,6] = 0
$data[62,7] = 28
$data[63,0] = "Argelia"
$data[63,1] = 53998
$data[63,2] = 221
$data[63,3] = 37856
$data[63,4] = 14301
$data[63,5] = 0
$data[63,6] = 14
$data[63,7] = 1841
$data[64,0] = "Paraguay"
$data[64,1] = 53482
$data[64,2] = 886
$data[64,3] = 34927
$data[64,4] = 17390
$data[64,5] = 0
$data[64,6] = 15
$data[64,7] = 1165
$data[65,0] = "Kirguistan"
$data[65,1] = 51020
$data[65,2] = 431
$data[65,3] = 45288
$data[65,4] = 4629
$data[65,5] = 0
$data[65,6] = 4
$data[65,7] = 1103
$data[66,0] = "Libia"
$data[66,1] = 47845
$data[66,2] = 1169
$data[66,3] = 26062
$data[66,4] = 21084
$data[66,5] = 0
$data[66,6] = 18
$data[66,7] = 699
$data[67,0] = "Irlanda"
$data[67,1] = 47427
$data[67,2] = 998
$data[67,3] = 23364
$data[67,4] = 22222
$data[67,5] = 0
$data[67,6] = 3
$data[67,7] = 1841
$data[68,0] = "Ghana"
$data[68,1] = 47173
$data[68,2] = 0
$data[68,3] = 46527
$data[68,4] = 336
$data[68,5] = 0
$data[68,6] = 0
$data[68,7] = 310
$data[69,0] = "Estado de Palestina"
$data[69,1] = 46434
$data[69,2] = 334
$data[69,3] = 39921
$data[69,4] = 6111
$data[69,5] = 0
$data[69,6] = 1
$data[69,7] = 402
$data[70,0] = "Azerbaiyan"
$data[70,1] = 43789
$data[70,2] = 509
$data[70,3] = 39800
$data[70,4] = 3368
$data[70,5] = 0
$data[70,6] = 2
$data[70,7] = 621
$data[71,0] = "Kenia"
$data[71,1] = 43580
$data[71,2] = 437
$data[71,3] = 31648
$data[71,4] = 11119
$data[71,5] = 0
$data[71,6] = 8
$data[71,7] = 813
$data[72,0] = "Hungria"
$data[72,1] = 43025
$data[72,2] = 1293
$data[72,3] = 13134
$data[72,4] = 28806
$data[72,5] = 0
$data[72,6] = 33
$data[72,7] = 1085
$data[73,0] = "Afganistan"
$data[73,1] = 40073
$data[73,2] = 47
$data[73,3] = 33516
$data[73,4] = 5072
$data[73,5] = 0
$data[73,6] = 4
$data[73,7] = 1485
$data[74,0] = "Serbia"
$data[74,1] = 35719
$data[74,2] = 265
$data[74,3] = 31536
$data[74,4] = 3411
$data[74,5] = 0
$data[74,6] = 2
$data[74,7] = 772
$data[75,0] = "Tunez"
$data[75,1] = 34790
$data[75,2] = 0
$data[75,3] = 5032
$data[75,4] = 29246
$data[75,5] = 0
$data[75,6] = 0
$data[75,7] = 512
$data[76,0] = "Jordania"
$data[76,1] = 34548
$data[76,2] = 1539
$data[76,3] = 6692
$data[76,4] = 27546
$data[76,5] = 0
$data[76,6] = 28
$data[76,7] = 310
$data[77,0] = "Dinamarca"
$data[77,1] = 34441
$data[77,2] = 418
$data[77,3] = 28551
$data[77,4] = 5213
$data[77,5] = 0
$data[77,6] = 0
$data[77,7] = 677
$data[78,0] = "Birmania"
$data[78,1] = 33488
$data[78,2] = 1137
$data[78,3] = 15477
$data[78,4] = 17212
$data[78,5] = 0
$data[78,6] = 34
$data[78,7] = 799
$data[79,0] = "Bosnia y Herzegovina"
$data[79,1] = 32845
$data[79,2] = 621
$data[79,3] = 24603
$data[79,4] = 7262
$data[79,5] = 0
$data[79,6] = 8
$data[79,7] = 980
$data[80,0] = "El Salvador"
$data[80,1] = 31265
$data[80,2] = 204
$data[80,3] = 26542
$data[80,4] = 3811
$data[80,5] = 0
$data[80,6] = 4
$data[80,7] = 912
$data[81,0] = "Bulgaria"
$data[81,1] = 28505
$data[81,2] = 998
$data[81,3] = 16875
$data[81,4] = 10672
$data[81,5] = 0
$data[81,6] = 14
$data[81,7] = 958
$data[82,0] = "Australia"
$data[82,1] = 27371
$data[82,2] = 9
$data[82,3] = 25062
$data[82,4] = 1405
$data[82,5] = 0
$data[82,6] = 0
$data[82,7] = 904
$data[83,0] = "Eslovaquia"
$data[83,1] = 26300
$data[83,2] = 2075
$data[83,3] = 7182
$data[83,4] = 19047
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 71
$data[84,0] = "Corea del Sur"
$data[84,1] = 25035
$data[84,2] = 47
$data[84,3] = 23180
$data[84,4] = 1414
$data[84,5] = 0
$data[84,6] = 2
$data[84,7] = 441
$data[85,0] = "Grecia"
$data[85,1] = 24450
$data[85,2] = 503
$data[85,3] = 9989
$data[85,4] = 13971
$data[85,5] = 0
$data[85,6] = 8
$data[85,7] = 490
$data[86,0] = "Croacia"
$data[86,1] = 23665
$data[86,2] = 1131
$data[86,3] = 19087
$data[86,4] = 4233
$data[86,5] = 0
$data[86,6] = 1
$data[86,7] = 345
$data[87,0] = "Republica de Macedonia"
$data[87,1] = 22607
$data[87,2] = 437
$data[87,3] = 16949
$data[87,4] = 4837
$data[87,5] = 0
$data[87,6] = 6
$data[87,7] = 821
$data[88,0] = "Camerun"
$data[88,1] = 21441
$data[88,2] = 0
$data[88,3] = 20117
$data[88,4] = 901
$data[88,5] = 0
$data[88,6] = 0
$data[88,7] = 423
$data[89,0] = "Costa de Marfil"
$data[89,1] = 20275
$data[89,2] = 18
$data[89,3] = 19953
$data[89,4] = 201
$data[89,5] = 0
$data[89,6] = 1
$data[89,7] = 121
$data[90,0] = "Malasia"
$data[90,1] = 18758
$data[90,2] = 629
$data[90,3] = 12259
$data[90,4] = 6323
$data[90,5] = 0
$data[90,6] = 6
$data[90,7] = 176
$data[91,0] = "Madagascar"
$data[91,1] = 16754
$data[91,2] = 0
$data[91,3] = 16124
$data[91,4] = 393
$data[91,5] = 0
$data[91,6] = 0
$data[91,7] = 237
$data[92,0] = "Albania"
$data[92,1] = 16501
$data[92,2] = 289
$data[92,3] = 9957
$data[92,4] = 6101
$data[92,5] = 0
$data[92,6] = 4
$data[92,7] = 443
$data[93,0] = "Noruega"
$data[93,1] = 16272
$data[93,2] = 136
$data[93,3] = 11863
$data[93,4] = 4131
$data[93,5] = 0
$data[93,6] = 0
$data[93,7] = 278
$data[94,0] = "Zambia"
$data[94,1] = 15659
$data[94,2] = 0
$data[94,3] = 14899
$data[94,4] = 414
$data[94,5] = 0
$data[94,6] = 0
$data[94,7] = 346
$data[95,0] = "Senegal"
$data[95,1] = 15368
$data[95,2] = 20
$data[95,3] = 13704
$data[95,4] = 1347
$data[95,5] = 0
$data[95,6] = 1
$data[95,7] = 317
$data[96,0] = "Georgia"
$data[96,1] = 15327
$data[96,2] = 887
$data[96,3] = 7613
$data[96,4] = 7590
$data[96,5] = 0
$data[96,6] = 11
$data[96,7] = 124
$data[97,0] = "Montenegro"
$data[97,1] = 15281
$data[97,2] = 273
$data[97,3] = 10569
$data[97,4] = 4484
$data[97,5] = 0
$data[97,6] = 7
$data[97,7] = 228
$data[98,0] = "Sudan"
$data[98,1] = 13691
$data[98,2] = 0
$data[98,3] = 6764
$data[98,4] = 6091
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 836
$data[99,0] = "Finlandia"
$data[99,1] = 13133
$data[99,2] = 189
$data[99,3] = 9100
$data[99,4] = 3682
$data[99,5] = 0
$data[99,6] = 1
$data[99,7] = 351
$data[100,0] = "Namibia"
$data[100,1] = 12215
$data[100,2] = 112
$data[100,3] = 10360
$data[100,4] = 1724
$data[100,5] = 0
$data[100,6] = 1
$data[100,7] = 131
$data[101,0] = "Eslovenia"
$data[101,1] = 11517
$data[101,2] = 834
$data[101,3] = 5924
$data[101,4] = 5413
$data[101,5] = 0
$data[101,6] = 4
$data[101,7] = 180
$data[102,0] = "Guinea"
$data[102,1] = 11362
$data[102,2] = 107
$data[102,3] = 10420
$data[102,4] = 872
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 70
$data[103,0] = "Maldivas"
$data[103,1] = 11154
$data[103,2] = 41
$data[103,3] = 9995
$data[103,4] = 1124
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 35
$data[104,0] = "Consejo Danes para los Refugiados"
$data[104,1] = 10999
$data[104,2] = 64
$data[104,3] = 10342
$data[104,4] = 356
$data[104,5] = 0
$data[104,6] = 20
$data[104,7] = 301
$data[105,0] = "Mozambique"
$data[105,1] = 10612
$data[105,2] = 75
$data[105,3] = 8262
$data[105,4] = 2277
$data[105,5] = 0
$data[105,6] = 0
$data[105,7] = 73
$data[106,0] = "Luxemburgo"
$data[106,1] = 10471
$data[106,2] = 227
$data[106,3] = 8468
$data[106,4] = 1870
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 133
$data[107,0] = "Tayikistan"
$data[107,1] = 10414
$data[107,2] = 40
$data[107,3] = 9393
$data[107,4] = 941
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 80
$data[108,0] = "Uganda"
$data[108,1] = 10334
$data[108,2] = 217
$data[108,3] = 6901
$data[108,4] = 3337
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 96
$data[109,0] = "Guayana Francesa"
$data[109,1] = 10239
$data[109,2] = 6
$data[109,3] = 9955
$data[109,4] = 215
$data[109,5] = 0
$data[109,6] = 0
$data[109,7] = 69
$data[110,0] = "Haiti"
$data[110,1] = 8925
$data[110,2] = 0
$data[110,3] = 7182
$data[110,4] = 1512
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 231
$data[111,0] = "Gabon"
$data[111,1] = 8881
$data[111,2] = 12
$data[111,3] = 8430
$data[111,4] = 397
$data[111,5] = 0
$data[111,6] = 0
$data[111,7] = 54
$data[112,0] = "Jamaica"
$data[112,1] = 8132
$data[112,2] = 65
$data[112,3] = 3653
$data[112,4] = 4317
$data[112,5] = 0
$data[112,6] = 2
$data[112,7] = 162
$data[113,0] = "Zimbabue"
$data[113,1] = 8099
$data[113,2] = 24
$data[113,3] = 7673
$data[113,4] = 195
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 231
$data[114,0] = "Mauritania"
$data[114,1] = 7603
$data[114,2] = 18
$data[114,3] = 7339
$data[114,4] = 101
$data[114,5] = 0
$data[114,6] = 0
$data[114,7] = 163
$data[115,0] = "Cabo Verde"
$data[115,1] = 7526
$data[115,2] = 82
$data[115,3] = 6425
$data[115,4] = 1019
$data[115,5] = 0
$data[115,6] = 3
$data[115,7] = 82
$data[116,0] = "Angola"
$data[116,1] = 7222
$data[116,2] = 126
$data[116,3] = 3012
$data[116,4] = 3976
$data[116,5] = 0
$data[116,6] = 6
$data[116,7] = 234
$data[117,0] = "Guadalupe"
$data[117,1] = 7122
$data[117,2] = 214
$data[117,3] = 2199
$data[117,4] = 4827
$data[117,5] = 0
$data[117,6] = 0
$data[117,7] = 96
$data[118,0] = "Lituania"
$data[118,1] = 7041
$data[118,2] = 281
$data[118,3] = 3035
$data[118,4] = 3894
$data[118,5] = 0
$data[118,6] = 2
$data[118,7] = 112
$data[119,0] = "Cuba"
$data[119,1] = 6118
$data[119,2] = 56
$data[119,3] = 5702
$data[119,4] = 292
$data[119,5] = 0
$data[119,6] = 1
$data[119,7] = 124
$data[120,0] = "Malaui"
$data[120,1] = 5842
$data[120,2] = 6
$data[120,3] = 4735
$data[120,4] = 926
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 181
$data[121,0] = "Suazilandia"
$data[121,1] = 5746
$data[121,2] = 13
$data[121,3] = 5392
$data[121,4] = 239
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 115
$data[122,0] = "Bahamas"
$data[122,1] = 5517
$data[122,2] = 132
$data[122,3] = 3201
$data[122,4] = 2202
$data[122,5] = 0
$data[122,6] = 2
$data[122,7] = 114
$data[123,0] = "Republica de Yibuti"
$data[123,1] = 5449
$data[123,2] = 6
$data[123,3] = 5372
$data[123,4] = 16
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 61
$data[124,0] = "Sri Lanka"
$data[124,1] = 5354
$data[124,2] = 110
$data[124,3] = 3385
$data[124,4] = 1956
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 13
$data[125,0] = "Nicaragua"
$data[125,1] = 5353
$data[125,2] = 0
$data[125,3] = 4225
$data[125,4] = 974
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 154
$data[126,0] = "Botsuana"
$data[126,1] = 5242
$data[126,2] = 0
$data[126,3] = 905
$data[126,4] = 4317
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 20
$data[127,0] = "Trinidad yTobago"
$data[127,1] = 5241
$data[127,2] = 47
$data[127,3] = 3545
$data[127,4] = 1601
$data[127,5] = 0
$data[127,6] = 2
$data[127,7] = 95
$data[128,0] = "Hong Kong"
$data[128,1] = 5221
$data[128,2] = 7
$data[128,3] = 4951
$data[128,4] = 165
$data[128,5] = 0
$data[128,6] = 0
$data[128,7] = 105
$data[129,0] = "Congo"
$data[129,1] = 5156
$data[129,2] = 0
$data[129,3] = 3887
$data[129,4] = 1177
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 92
$data[130,0] = "Surinam"
$data[130,1] = 5113
$data[130,2] = 19
$data[130,3] = 4921
$data[130,4] = 83
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 109
$data[131,0] = "Guinea Ecuatorial"
$data[131,1] = 5068
$data[131,2] = 0
$data[131,3] = 4954
$data[131,4] = 31
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 83
$data[132,0] = "Siria"
$data[132,1] = 4987
$data[132,2] = 56
$data[132,3] = 1456
$data[132,4] = 3290
$data[132,5] = 0
$data[132,6] = 3
$data[132,7] = 241
$data[133,0] = "Ruanda"
$data[133,1] = 4965
$data[133,2] = 12
$data[133,3] = 4664
$data[133,4] = 267
$data[133,5] = 0
$data[133,6] = 1
$data[133,7] = 34
$data[134,0] = "Republica de Africa Central"
$data[134,1] = 4855
$data[134,2] = 0
$data[134,3] = 1924
$data[134,4] = 2869
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 62
$data[135,0] = "Reunion"
$data[135,1] = 4776
$data[135,2] = 98
$data[135,3] = 4445
$data[135,4] = 314
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 17
$data[136,0] = "Aruba"
$data[136,1] = 4289
$data[136,2] = 4
$data[136,3] = 3947
$data[136,4] = 310
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 32
$data[137,0] = "Malta"
$data[137,1] = 4282
$data[137,2] = 122
$data[137,3] = 3142
$data[137,4] = 1095
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 45
$data[138,0] = "Mayotte"
$data[138,1] = 4030
$data[138,2] = 0
$data[138,3] = 2964
$data[138,4] = 1023
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 43
$data[139,0] = "Estonia"
$data[139,1] = 4017
$data[139,2] = 37
$data[139,3] = 3137
$data[139,4] = 812
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 68
$data[140,0] = "Islandia"
$data[140,1] = 3929
$data[140,2] = 92
$data[140,3] = 2713
$data[140,4] = 1205
$data[140,5] = 0
$data[140,6] = 1
$data[140,7] = 11
$data[141,0] = "Somalia"
$data[141,1] = 3864
$data[141,2] = 0
$data[141,3] = 3089
$data[141,4] = 676
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 99
$data[142,0] = "Polinesia Francesa"
$data[142,1] = 3797
$data[142,2] = 224
$data[142,3] = 2844
$data[142,4] = 939
$data[142,5] = 0
$data[142,6] = 1
$data[142,7] = 14
$data[143,0] = "Guyana"
$data[143,1] = 3672
$data[143,2] = 52
$data[143,3] = 2590
$data[143,4] = 975
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 107
$data[144,0] = "Tailandia"
$data[144,1] = 3669
$data[144,2] = 4
$data[144,3] = 3467
$data[144,4] = 143
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 59
$data[145,0] = "Gambia"
$data[145,1] = 3649
$data[145,2] = 5
$data[145,3] = 2649
$data[145,4] = 882
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 118
$data[146,0] = "Mali"
$data[146,1] = 3378
$data[146,2] = 10
$data[146,3] = 2563
$data[146,4] = 683
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 132
$data[147,0] = "Principado de Andorra"
$data[147,1] = 3377
$data[147,2] = 187
$data[147,3] = 2057
$data[147,4] = 1261
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 59
$data[148,0] = "Letonia"
$data[148,1] = 3204
$data[148,2] = 148
$data[148,3] = 1329
$data[148,4] = 1833
$data[148,5] = 0
$data[148,6] = 1
$data[148,7] = 42
$data[149,0] = "Sudan del Sur"
$data[149,1] = 2817
$data[149,2] = 10
$data[149,3] = 1290
$data[149,4] = 1472
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 55
$data[150,0] = "Belice"
$data[150,1] = 2682
$data[150,2] = 63
$data[150,3] = 1612
$data[150,4] = 1029
$data[150,5] = 0
$data[150,6] = 1
$data[150,7] = 41
$data[151,0] = "Benin"
$data[151,1] = 2496
$data[151,2] = 0
$data[151,3] = 2330
$data[151,4] = 125
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 41
$data[152,0] = "Uruguay"
$data[152,1] = 2450
$data[152,2] = 33
$data[152,3] = 2042
$data[152,4] = 357
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 51
$data[153,0] = "Guinea-Bisau"
$data[153,1] = 2389
$data[153,2] = 0
$data[153,3] = 1782
$data[153,4] = 566
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 41
$data[154,0] = "Republica de Chipre"
$data[154,1] = 2379
$data[154,2] = 94
$data[154,3] = 1444
$data[154,4] = 910
$data[154,5] = 0
$data[154,6] = 0
$data[154,7] = 25
$data[155,0] = "Burkina Faso"
$data[155,1] = 2343
$data[155,2] = 8
$data[155,3] = 1718
$data[155,4] = 560
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 65
$data[156,0] = "Sierra Leona"
$data[156,1] = 2325
$data[156,2] = 2
$data[156,3] = 1750
$data[156,4] = 502
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 73
$data[157,0] = "Martinica"
$data[157,1] = 2257
$data[157,2] = 0
$data[157,3] = 98
$data[157,4] = 2135
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 24
$data[158,0] = "Yemen"
$data[158,1] = 2055
$data[158,2] = 2
$data[158,3] = 1335
$data[158,4] = 124
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 596
$data[159,0] = "Togo"
$data[159,1] = 2027
$data[159,2] = 31
$data[159,3] = 1500
$data[159,4] = 476
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 51
$data[160,0] = "Nueva Zelanda"
$data[160,1] = 1880
$data[160,2] = 4
$data[160,3] = 1809
$data[160,4] = 46
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 25
$data[161,0] = "Lesoto"
$data[161,1] = 1833
$data[161,2] = 0
$data[161,3] = 961
$data[161,4] = 830
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 42
$data[162,0] = "Liberia"
$data[162,1] = 1377
$data[162,2] = 3
$data[162,3] = 1264
$data[162,4] = 31
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 82
$data[163,0] = "Republica del Chad"
$data[163,1] = 1361
$data[163,2] = 11
$data[163,3] = 1138
$data[163,4] = 130
$data[163,5] = 0
$data[163,6] = 1
$data[163,7] = 93
$data[164,0] = "Niger"
$data[164,1] = 1209
$data[164,2] = 2
$data[164,3] = 1126
$data[164,4] = 14
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 69
$data[165,0] = "Vietnam"
$data[165,1] = 1124
$data[165,2] = 0
$data[165,3] = 1031
$data[165,4] = 58
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 35
$data[166,0] = "Santo Tome y Principe"
$data[166,1] = 932
$data[166,2] = 3
$data[166,3] = 896
$data[166,4] = 21
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 15
$data[167,0] = "San Marino"
$data[167,1] = 759
$data[167,2] = 0
$data[167,3] = 685
$data[167,4] = 32
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 42
$data[168,0] = "San Martin (Parte Holandesa)"
$data[168,1] = 746
$data[168,2] = 9
$data[168,3] = 659
$data[168,4] = 65
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 22
$data[169,0] = "Crucero"
$data[169,1] = 712
$data[169,2] = 0
$data[169,3] = 659
$data[169,4] = 40
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 13
$data[170,0] = "Curazao"
$data[170,1] = 698
$data[170,2] = 25
$data[170,3] = 391
$data[170,4] = 306
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 1
$data[171,0] = "Islas Turcas y Caicos"
$data[171,1] = 697
$data[171,2] = 1
$data[171,3] = 674
$data[171,4] = 17
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 6
$data[172,0] = "Papua Nueva Guinea"
$data[172,1] = 578
$data[172,2] = 0
$data[172,3] = 537
$data[172,4] = 34
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 7
$data[173,0] = "Gibraltar"
$data[173,1] = 544
$data[173,2] = 13
$data[173,3] = 441
$data[173,4] = 103
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 0
$data[174,0] = "Taiwan"
$data[174,1] = 535
$data[174,2] = 4
$data[174,3] = 491
$data[174,4] = 37
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 7
$data[175,0] = "San Martin (Parte Francesa)"
$data[175,1] = 531
$data[175,2] = 30
$data[175,3] = 380
$data[175,4] = 143
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 8
$data[176,0] = "Burundi"
$data[176,1] = 531
$data[176,2] = 2
$data[176,3] = 497
$data[176,4] = 33
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 1
$data[177,0] = "Tanzania"
$data[177,1] = 509
$data[177,2] = 0
$data[177,3] = 183
$data[177,4] = 305
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 21
$data[178,0] = "Comoras"
$data[178,1] = 502
$data[178,2] = 6
$data[178,3] = 485
$data[178,4] = 10
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 7
$data[179,0] = "Islas Feroe"
$data[179,1] = 482
$data[179,2] = 2
$data[179,3] = 471
$data[179,4] = 11
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 0
$data[180,0] = "Eritrea"
$data[180,1] = 422
$data[180,2] = 0
$data[180,3] = 376
$data[180,4] = 46
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 0
$data[181,0] = "Mauricio"
$data[181,1] = 417
$data[181,2] = 2
$data[181,3] = 364
$data[181,4] = 43
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 10
$data[182,0] = "Isla de Man"
$data[182,1] = 348
$data[182,2] = 0
$data[182,3] = 319
$data[182,4] = 5
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 24
$data[183,0] = "Mongolia"
$data[183,1] = 320
$data[183,2] = 0
$data[183,3] = 311
$data[183,4] = 9
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = "Butan"
$data[184,1] = 316
$data[184,2] = 0
$data[184,3] = 296
$data[184,4] = 20
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = "Camboya"
$data[185,1] = 283
$data[185,2] = 0
$data[185,3] = 278
$data[185,4] = 5
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 0
$data[186,0] = "Monaco"
$data[186,1] = 255
$data[186,2] = 2
$data[186,3] = 217
$data[186,4] = 36
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 2
$data[187,0] = "Islas Caimanes"
$data[187,1] = 233
$data[187,2] = 8
$data[187,3] = 212
$data[187,4] = 20
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 1
$data[188,0] = "Barbados"
$data[188,1] = 219
$data[188,2] = 1
$data[188,3] = 195
$data[188,4] = 17
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 7
$data[189,0] = "Liechtenstein"
$data[189,1] = 192
$data[189,2] = 9
$data[189,3] = 132
$data[189,4] = 59
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 1
$data[190,0] = "Bermudas"
$data[190,1] = 185
$data[190,2] = 0
$data[190,3] = 172
$data[190,4] = 4
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 9
$data[191,0] = "Bonaire, San Eustaquio y Saba"
$data[191,1] = 150
$data[191,2] = 0
$data[191,3] = 111
$data[191,4] = 37
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 2
$data[192,0] = "Seychelles"
$data[192,1] = 149
$data[192,2] = 1
$data[192,3] = 148
$data[192,4] = 1
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 0
$data[193,0] = "Brunei"
$data[193,1] = 147
$data[193,2] = 0
$data[193,3] = 143
$data[193,4] = 1
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 3
$data[194,0] = "Antigua y Barbuda"
$data[194,1] = 112
$data[194,2] = 0
$data[194,3] = 100
$data[194,4] = 9
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 3
$data[195,0] = "San Bartolome"
$data[195,1] = 72
$data[195,2] = 5
$data[195,3] = 55
$data[195,4] = 17
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 0
$data[196,0] = "Islas Virgenes Britanicas"
$data[196,1] = 71
$data[196,2] = 0
$data[196,3] = 70
$data[196,4] = 0
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 1
$data[197,0] = "San Vicente y las Granadinas"
$data[197,1] = 65
$data[197,2] = 1
$data[197,3] = 64
$data[197,4] = 1
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = "Macao"
$data[198,1] = 46
$data[198,2] = 0
$data[198,3] = 46
$data[198,4] = 0
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 0
$data[199,0] = "Puerto Rico"
$data[199,1] = 39
$data[199,2] = 0
$data[199,3] = 1
$data[199,4] = 36
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 2
$data[200,0] = "Dominica"
$data[200,1] = 33
$data[200,2] = 0
$data[200,3] = 29
$data[200,4] = 4
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 0
$data[201,0] = "Guam"
$data[201,1] = 32
$data[201,2] = 0
$data[201,3] = 0
$data[201,4] = 31
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 1
$data[202,0] = "Santa Lucia"
$data[202,1] = 32
$data[202,2] = 1
$data[202,3] = 27
$data[202,4] = 5
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = "Fiyi"
$data[203,1] = 32
$data[203,2] = 0
$data[203,3] = 30
$data[203,4] = 0
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 2
$data[204,0] = "Timor Oriental"
$data[204,1] = 29
$data[204,2] = 0
$data[204,3] = 28
$data[204,4] = 1
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = "Nueva Caledonia"
$data[205,1] = 27
$data[205,2] = 0
$data[205,3] = 27
$data[205,4] = 0
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = "Santa Sede"
$data[206,1] = 26
$data[206,2] = 0
$data[206,3] = 12
$data[206,4] = 14
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 0
$data[207,0] = "Granada"
$data[207,1] = 25
$data[207,2] = 0
$data[207,3] = 24
$data[207,4] = 1
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 0
$data[208,0] = "Laos"
$data[208,1] = 23
$data[208,2] = 0
$data[208,3] = 22
$data[208,4] = 1
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 0
$data[209,0] = "San Cristobal y Nieves"
$data[209,1] = 19
$data[209,2] = 0
$data[209,3] = 19
$data[209,4] = 0
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 0
$data[210,0] = "Islas Virgenes de los Estados Unidos"
$data[210,1] = 17
$data[210,2] = 0
$data[210,3] = 0
$data[210,4] = 17
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 0
$data[211,0] = "San Pedro y Miquelon"
$data[211,1] = 16
$data[211,2] = 0
$data[211,3] = 12
$data[211,4] = 4
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 0
$data[212,0] = "Groenlandia"
$data[212,1] = 16
$data[212,2] = 0
$data[212,3] = 14
$data[212,4] = 2
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 0
$data[213,0] = "Montserrat"
$data[213,1] = 13
$data[213,2] = 0
$data[213,3] = 12
$data[213,4] = 0
$data[213,5] = 0
$data[213,6] = 0
$data[213,7] = 1
$data[214,0] = "Islas Malvinas"
$data[214,1] = 13
$data[214,2] = 0
$data[214,3] = 13
$data[214,4] = 0
$data[214,5] = 0
$data[214,6] = 0
$data[214,7] = 0
$data[215,0] = "Sahara Occidental"
$data[215,1] = 10
$data[215,2] = 0
$data[215,3] = 8
$data[215,4] = 1
$data[215,5] = 0
$data[215,6] = 0
$data[215,7] = 1
$data[216,0] = "Islas Salomon"
$data[216,1] = 3
$data[216,2] = 0
$data[216,3] = 0
$data[216,4] = 3
$data[216,5] = 0
$data[216,6] = 0
$data[216,7] = 0
$data[217,0] = "Anguila"
$data[217,1] = 3
$data[217,2] = 0
$data[217,3] = 3
$data[217,4] = 0
$data[217,5] = 0
$data[217,6] = 0
$data[217,7] = 0
$data[218,0] = "Wallis y Futuna"
$data[218,1] = 1
$data[218,2] = 1
$data[218,3] = 0
$data[218,4] = 1
$data[218,5] = 0
$data[218,6] = 0
$data[218,7] = 0

$ws.Range("A3:H221").Value = $data

